# Apply the Metadata sheet update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to 2022-01-21T20:46:54+00:00
#  - Publisher/Contact rows replaced with Publisher/Jurisdiction rows
#    (the duplicated "Contact"/"No display for ContactDetail" row is removed,
#     and a new "Jurisdiction" / "United States of America" row takes its place)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (row 11); this shifts rows 12-15 up to 11-14
$ws.Rows.Item(11).Delete()

# Update Version value
$ws.Range("B3").Value = "6.0.0"

# Update Date value
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Fill in Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Replace the remaining Contact row with the new Jurisdiction row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
